# Adds a new "2022-Q4" quarter sheet (with fund holding detail data) right
# after the "总计" (totals) sheet, and updates the "总计" sheet with a new
# summary row for 2022-Q4 (shifting the existing quarters' summary rows
# down by one data slot, so the oldest quarter - 2021-Q1 - gets its own
# new trailing row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet, positioned right after "总计"
#    (i.e. right before the current "2022-Q3" sheet).
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3Sheet)
$newSheet.Name = "2022-Q4"

# Header row for the new sheet.
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Copy the header style (bold + centered + boxed, style index "2" in the
# original file) from the "总计" sheet's own header so formatting matches.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Fund holding detail rows for 2022-Q4.
$q4Data = @(
  @(0,  "090018", "大成新锐产业混合",             "88.75", "93.33", "4.11", "3.6476", 7),
  @(1,  "001300", "大成睿景灵活配置混合A",         "32.95", "92.29", "4.31", "1.4201", 7),
  @(2,  "013435", "大成景气精选六个月持有混合A",   "30.45", "91.16", "3.15", "0.9592", 9),
  @(3,  "001301", "大成睿景灵活配置混合C",         "19.30", "92.29", "4.31", "0.8318", 7),
  @(4,  "002258", "大成国企改革灵活配置混合",       "16.71", "93.37", "4.22", "0.7052", 7),
  @(5,  "014224", "大成聚优成长混合A",             "16.55", "90.21", "3.18", "0.5263", 9),
  @(6,  "010826", "大成产业趋势混合A",             "11.37", "93.99", "4.43", "0.5037", 6),
  @(7,  "012519", "大成核心趋势混合A",             "9.51",  "91.16", "3.13", "0.2977", 9),
  @(8,  "013436", "大成景气精选六个月持有混合C",   "5.20",  "91.16", "3.15", "0.1638", 9),
  @(9,  "010827", "大成产业趋势混合C",             "3.42",  "93.99", "4.43", "0.1515", 6),
  @(10, "014225", "大成聚优成长混合C",             "3.26",  "90.21", "3.18", "0.1037", 9),
  @(11, "002295", "广发稳安灵活配置混合A",         "1.51",  "80.18", "5.49", "0.0829", 2),
  @(12, "012520", "大成核心趋势混合C",             "2.45",  "91.16", "3.13", "0.0767", 9),
  @(13, "005770", "信澳中证沪港深高股息精选指数",   "0.13",  "23.47", "0.51", "0.0007", 10),
  @(14, "008604", "广发稳安灵活配置混合C",         "0.01",  "80.18", "5.49", "0.0005", 2)
)

$row = 2
foreach ($rec in $q4Data) {
    $newSheet.Cells.Item($row, 1).Value = [double]$rec[0]
    $newSheet.Cells.Item($row, 2).Value = "'" + $rec[1]
    $newSheet.Cells.Item($row, 3).Value = $rec[2]
    $newSheet.Cells.Item($row, 4).Value = "'" + $rec[3]
    $newSheet.Cells.Item($row, 5).Value = "'" + $rec[4]
    $newSheet.Cells.Item($row, 6).Value = "'" + $rec[5]
    $newSheet.Cells.Item($row, 7).Value = "'" + $rec[6]
    $newSheet.Cells.Item($row, 8).Value = [double]$rec[7]

    # Match the "A" (index) column's header-like style used throughout
    # the rest of the workbook for that column.
    $totalSheet.Range("A2").Copy()
    $newSheet.Cells.Item($row, 1).PasteSpecial(-4122)

    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) summary sheet: shift the quarterly rows
#    down (B/C/D columns) so a new 2022-Q4 entry lands on row 2, and the
#    previously-last quarter (2021-Q1) gets appended as a brand-new row.
# ---------------------------------------------------------------------
$quarters = @(
  @("2022-Q4", 15, 9.47),
  @("2022-Q3", 14, 8.21),
  @("2022-Q2", 6,  0.86),
  @("2022-Q1", 8,  9.21),
  @("2021-Q4", 12, 12.31),
  @("2021-Q3", 18, 14.45),
  @("2021-Q2", 22, 10.59),
  @("2021-Q1", 14, 4.83)
)

# Style for the new A9 index cell should match the rest of column A.
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

$r = 2
foreach ($q in $quarters) {
    $totalSheet.Cells.Item($r, 2).Value = $q[0]
    $totalSheet.Cells.Item($r, 3).Value = $q[1]
    $totalSheet.Cells.Item($r, 4).Value = $q[2]
    $r = $r + 1
}
